$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Demonstrate the different ways a boolean ("batch_useable", column H) can be
# authored in the sheet: true/false, 0/1, and yes/no/y/n strings.
# Row 6  (H6)  stays TRUE/FALSE boolean -> unchanged
$ws.Range("H7").Value = 0
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = "Yes"
$ws.Range("H10").Value = "No"
$ws.Range("H11").Value = "Y"
$ws.Range("H12").Value = "N"
# Rows 13-15 (H13, H14, H15) stay TRUE/FALSE boolean -> unchanged

# Move the active selection the way the author left it after editing.
$ws.Range("H7").Select()
